$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper constant: xlPasteValues
$xlPasteValues = -4163

function Set-TextValue($range, [string]$text) {
    # Assigning a plain string to .Value lets Excel auto-detect dates/numbers
    # (e.g. "01/01/2023" becomes a date serial). Route the literal text
    # through a quoted formula first, then collapse it back down to a
    # static value via copy / paste-special-values so the cell keeps its
    # existing style and is stored as literal text.
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

# Row 10 ("Objetivos:") — now holds the first responsible teacher
Set-TextValue $ws.Range("B10") "519033 - Carlos Yujiro Shigue"
Set-TextValue $ws.Range("C10") "519033 - Carlos Yujiro Shigue"

# Row 13 ("Programa resumido:") — now holds the activation date
Set-TextValue $ws.Range("B13") "01/01/2023"
Set-TextValue $ws.Range("C13") "01/01/2023"

# Row 15 ("Programa:") — now holds the first responsible teacher (duplicate)
Set-TextValue $ws.Range("B15") "519033 - Carlos Yujiro Shigue"
Set-TextValue $ws.Range("C15") "519033 - Carlos Yujiro Shigue"

# Row 18 ("Método:") — now holds the second responsible teacher
Set-TextValue $ws.Range("B18") "7290967 - Emerson Gonçalves de Melo"
Set-TextValue $ws.Range("C18") "7290967 - Emerson Gonçalves de Melo"
